$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "asdas"
$ws.Range("C2").Value = "asd"
$ws.Range("D2").Value = "asd"
$ws.Range("E2").Value = "asd"
$ws.Range("F2").Value = "sd"
$ws.Range("G2").Value = "sd"
$ws.Range("H2").Value = "sad"
$ws.Range("I2").Value = "sd"
$ws.Range("I2").NumberFormat = '_-"$"* #,##0.00_-;\-"$"* #,##0.00_-;_-"$"* "-"??_-;_-@_-'

$ws.Range("I2").Select()
